{"js": "// Append three new runs of text to the end of the second paragraph\n// (the \"Systemet samler...\" / \"...forretningsrejser.\" paragraph),\n// right after the existing text, matching the author's edit:\n//   \" P\u00e5 baggrund af en it-forunders\u00f8gelse \"\n//   \"kan vi se at processen vil kunne blive optimeret v\u00e6sentligt med et it system. Vi vil kunne udelukke flere folk fra processen og dermed optimere markant.\"\n//   \" \"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that ends with the known sentence about\n// \"salgschefens forretningsrejser.\" so the edit is anchored on content,\n// not a hard-coded index.\nconst marker = \"salgschefens forretningsrejser.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find target paragraph containing: \" + marker);\n}\n\ntarget.insertText(\n  \" P\u00e5 baggrund af en it-forunders\u00f8gelse \",\n  Word.InsertLocation.end\n);\ntarget.insertText(\n  \"kan vi se at processen vil kunne blive optimeret v\u00e6sentligt med et it system. Vi vil kunne udelukke flere folk fra processen og dermed optimere markant.\",\n  Word.InsertLocation.end\n);\ntarget.insertText(\" \", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Append three new runs of text to the end of the paragraph that ends\n# with \"...salgschefens forretningsrejser.\", matching the author's edit:\n#   \" P\u00e5 baggrund af en it-forunders\u00f8gelse \"\n#   \"kan vi se at processen vil kunne blive optimeret v\u00e6sentligt med et it system. Vi vil kunne udelukke flere folk fra processen og dermed optimere markant.\"\n#   \" \"\n$d = $word.ActiveDocument\n\n$marker = \"salgschefens forretningsrejser.\"\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like (\"*\" + $marker + \"*\")) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find target paragraph containing: $marker\"\n}\n\n# Re-fetch the paragraph's Range before each insertion so the Range\n# reflects the paragraph's (growing) extent and the new text lands\n# inside the same paragraph, before its paragraph mark.\n$p1 = $d.Paragraphs.Item($targetIndex)\n$p1.Range.InsertAfter(\" P\u00e5 baggrund af en it-forunders\u00f8gelse \")\n\n$p2 = $d.Paragraphs.Item($targetIndex)\n$p2.Range.InsertAfter(\"kan vi se at processen vil kunne blive optimeret v\u00e6sentligt med et it system. Vi vil kunne udelukke flere folk fra processen og dermed optimere markant.\")\n\n$p3 = $d.Paragraphs.Item($targetIndex)\n$p3.Range.InsertAfter(\" \")\n"}
